# UIUX.U39_CustomRetailUsers_action_init.xlsx
# Re-purpose the test-case sheet from the old "OCBOut_getCardDetail" SOAP
# test to the new "CustomRetailUsers_action=init" REST/OData test:
#   - Document Code / Function name
#   - Developer + System (IIB -> OCB)
#   - Input data (SOAP request -> REST URL)
#   - Output data (SOAP response -> JSON response)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$jsonResponse = @'
{
    "d": {
        "__metadata": {
            "id": "https://smp-srv:8081/cb/odata/services/retailuserservice/CustomRetailUsers('469')",
            "uri": "https://smp-srv:8081/cb/odata/services/retailuserservice/CustomRetailUsers('469')",
            "type": "com.sap.banking.custom.user.endpoint.v1_0.beans.CustomRetailUser"
        },
        "Id": "469",
        "FirstName": "none",
        "LastName": "VPBANK469",
        "MiddleName": null,
        "AffiliateBankID": 1,
        "BankId": "1000",
        "ConfirmPassword": null,
        "ConfirmPasswordReminder": null,
        "ConfirmPasswordReminder2": null,
        "CustId": "884206",
        "Greeting": null,
        "GreetingType": "1",
        "Password": "FFIHASHKXsZIfopn7b7u3sxzB6LbWJ5HLZfMhmAE+UhiJcCZKo=",
        "PasswordClue": "What is your favorite food?",
        "PasswordClue2": "What is your favorite movie?",
        "MemberId": null,
        "GroupId": "2139",
        "GroupName": null,
        "PasswordReminder": null,
        "PasswordReminder2": null,
        "PersonalBanker": "0",
        "Ssn": null,
        "MaskSSN": "",
        "UserName": "namnguyen",
        "AccountStatus": "1",
        "Timeout": "300",
        "CustomerType": "2",
        "RequestedCarrierTCId": null,
        "Address": {
            "__metadata": {
                "type": "com.sap.banking.common.endpoint.v1_0.beans.Address"
            },
            "Street": "VPBANK469",
            "Street2": "VPBANK469",
            "Street3": null,
            "City": null,
            "State": null,
            "StreetCode": null,
            "Country": "VNM",
            "Email": "tienthanhle38@gmail.com",
            "Phone": "03495896391",
            "Phone2": "VPBANK469",
            "ZipCode": null,
            "DataPhone": "0383962087",
            "FaxPhone": null,
            "PreferredContactMethod": null,
            "PreferredLanguage": "en_US"
        },
        "OldPassword": null,
        "NewPassword": null,
        "OldOtpMethod": "2",
        "NewOtpMethod": null,
        "OldServicePackage": "200000012",
        "OldServicePackageName": "PLATINIUM",
        "NewServicePackage": null,
        "NewServicePackageName": null,
        "ServicePackages": "200000010#INQUIRY;200000011#STANDARD",
        "Transactions": {
            "__deferred": {
                "uri": "https://smp-srv:8081/cb/odata/services/retailuserservice/CustomRetailUsers('469')/Transactions"
            }
        }
    }
}
'@

$requestUrl = "https://smp-srv:8081/cb/odata/services/retailuserservice/CustomRetailUsers('469')?action=init"

# Developer
$ws.Range("C9").Value = "khoand"

# System: IIB -> OCB (both the header summary and the table row)
$ws.Range("C10").Value = "OCB"
$ws.Range("D15").Value = "OCB"

# Document Code (header block) + the matching row in the test-case table
$ws.Range("C6").Value = "CustomRetailUsers_action=init"
$ws.Range("B15").Value = "CustomRetailUsers_action=init"

# Output data / Input data for the single test case row
$ws.Range("G15").Value = $jsonResponse
$ws.Range("F15").Value = $requestUrl

# Refresh the selection to match the edited row (author was working on G15)
$ws.Range("G15").Select()
